$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 150
$ws.Range("I21").Value = 150
$ws.Range("K21").Value = 150
$ws.Range("M21").Value = 318

$ws.Range("H23").Value = 150
$ws.Range("I23").Value = 150
$ws.Range("K23").Value = 150
$ws.Range("M23").Value = 84

$ws.Range("H111").Value = 4506
$ws.Range("I111").Value = 1029
$ws.Range("J111").Value = 5665
$ws.Range("K111").Value = 3087
$ws.Range("L111").Value = 16995
$ws.Range("M111").Value = -20
$ws.Range("N111").Value = -23129

$ws.Range("H116").Value = 3198.4
$ws.Range("J116").Value = 3992
$ws.Range("L116").Value = 3992
$ws.Range("N116").Value = -10876

$ws.Range("H125").Value = 5443.3335
$ws.Range("I125").Value = 1330
$ws.Range("K125").Value = 11970
$ws.Range("M125").Value = -9510

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1754.4
$ws.Range("I61").Value = 1754.4
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1754.4
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1542.4
$ws.Range("N61").ClearContents()

$ws.Range("H102").Value = 930.8
$ws.Range("I102").Value = 930.8
$ws.Range("K102").Value = 930.8
$ws.Range("M102").Value = 691.2

$ws.Range("H110").Value = 700
$ws.Range("I110").Value = 550
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 550
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1495
$ws.Range("N110").Value = -5090

$ws.Range("H136").Value = 1754.4
$ws.Range("I136").Value = 1754.4
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5263.200000000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2713.200000000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2060.625
$ws.Range("I20").Value = 1997.8
$ws.Range("K20").Value = 1997.8
$ws.Range("M20").Value = -1750.8

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H105").Value = 2599.4
$ws.Range("I105").Value = 2332.3333
$ws.Range("K105").Value = 2332.3333
$ws.Range("M105").Value = -585.3332999999998

$ws.Range("H134").Value = 3753
$ws.Range("I134").Value = 4528.6
$ws.Range("K134").Value = 13585.8
$ws.Range("M134").Value = -11050.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2050
$ws.Range("I31").Value = 700
$ws.Range("J31").Value = 3400
$ws.Range("K31").Value = 700
$ws.Range("L31").Value = 3400
$ws.Range("M31").Value = -405
$ws.Range("N31").Value = -3990

$ws.Range("H34").Value = 2050
$ws.Range("I34").Value = 700
$ws.Range("J34").Value = 3400
$ws.Range("K34").Value = 700
$ws.Range("L34").Value = 3400
$ws.Range("M34").Value = -498
$ws.Range("N34").Value = -3804

$ws.Range("H86").Value = 10829.333
$ws.Range("I86").Value = 9994.5
$ws.Range("K86").Value = 9994.5
$ws.Range("M86").Value = -8871.5

$ws.Range("H89").Value = 10829.333
$ws.Range("I89").Value = 9994.5
$ws.Range("K89").Value = 49972.5
$ws.Range("M89").Value = -44356.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1011.2
$ws.Range("I140").Value = 1011.2
$ws.Range("K140").Value = 3033.6
$ws.Range("M140").Value = 2146.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 2999.5
$ws.Range("I41").Value = 999
$ws.Range("K41").Value = 999
$ws.Range("M41").Value = -644

$ws.Range("H126").Value = 8003.5
$ws.Range("I126").Value = 8003.5
$ws.Range("K126").Value = 24010.5
$ws.Range("M126").Value = -21540.5

$ws.Range("H132").Value = 8004.6
$ws.Range("I132").Value = 8004.6
$ws.Range("K132").Value = 24013.8
$ws.Range("M132").Value = -21483.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4004
$ws.Range("I7").Value = 4004
$ws.Range("K7").Value = 4004
$ws.Range("M7").Value = -3892

$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10452

$ws.Range("H22").Value = 2620
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -255
$ws.Range("N22").Value = -4590

$ws.Range("H27").Value = 2620
$ws.Range("I27").Value = 550
$ws.Range("J27").Value = 4000
$ws.Range("K27").Value = 550
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = -443
$ws.Range("N27").Value = -4214

$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 3000
$ws.Range("L42").Value = 3000
$ws.Range("N42").Value = -4126

$ws.Range("H49").Value = 3000
$ws.Range("J49").Value = 3000
$ws.Range("L49").Value = 3000
$ws.Range("N49").Value = -3294

$ws.Range("H55").Value = 2122
$ws.Range("I55").Value = 999
$ws.Range("J55").Value = 2496.3333
$ws.Range("K55").Value = 999
$ws.Range("L55").Value = 2496.3333
$ws.Range("M55").Value = -826
$ws.Range("N55").Value = -2842.3333

$ws.Range("H126").Value = 4004
$ws.Range("I126").Value = 4004
$ws.Range("K126").Value = 12012
$ws.Range("M126").Value = -9542

$ws.Range("H132").Value = 8833.333000000001
$ws.Range("I132").Value = 14500
$ws.Range("K132").Value = 43500
$ws.Range("M132").Value = -40970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1499
$ws.Range("I14").Value = 1499
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1499
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1331
$ws.Range("N14").ClearContents()

$ws.Range("H41").Value = 64767.8
$ws.Range("I41").Value = 64638
$ws.Range("K41").Value = 64638
$ws.Range("M41").Value = -64248

$ws.Range("H126").Value = 1340.7273
$ws.Range("I126").Value = 1074.8
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 3224.4
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -754.3999999999996
$ws.Range("N126").Value = -16940

$ws.Range("H132").Value = 3907.25
$ws.Range("I132").Value = 3394
$ws.Range("K132").Value = 10182
$ws.Range("M132").Value = -7652
